{"js": "// Word Javascript API (Office.js) script.\n// Applies the same edits as the authoritative OOXML diff:\n//  1. \"happy-test\" -> \"happy path test\" in the \"I added one ...\" paragraph.\n//  2. A new paragraph about seeding the \"database\" is inserted right after\n//     that paragraph.\n//  3. The \"...in which case\" paragraph gets a new trailing sentence.\n//  4. The closing sentence of the document is expanded with \"behavioral \"\n//     and a trailing example clause.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) \"happy-test\" -> \"happy path test\"\n// ---------------------------------------------------------------------\nconst happyTestResults = body.search(\"happy-test\", { matchCase: true });\nhappyTestResults.load(\"text\");\nawait context.sync();\n\nif (happyTestResults.items.length > 0) {\n  happyTestResults.items[0].insertText(\"happy path test\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) Insert a new paragraph right after \"...case description.\" (the\n//    paragraph that starts with \"I added one happy path test ...\").\n// ---------------------------------------------------------------------\nconst sanityCheckResults = body.search(\n  \"I added one happy path test for each service function in order as a form of sanity check that I implemented all operations specified in the case description.\",\n  { matchCase: true }\n);\nsanityCheckResults.load(\"text\");\nawait context.sync();\n\nif (sanityCheckResults.items.length > 0) {\n  sanityCheckResults.items[0].insertParagraph(\n    \"To seed the \\u201cdatabase\\u201d I used the service layer calls, normally I would Mock or Fake this layer \" +\n      \"to seed data necessary for each test. I skipped that due to time restraint.\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) \"...in which case\" -> \"...in which case the actual implementation\n//    shall already be tested.\"\n// ---------------------------------------------------------------------\nconst inWhichCaseResults = body.search(\"in which case\", { matchCase: true });\ninWhichCaseResults.load(\"text\");\nawait context.sync();\n\nif (inWhichCaseResults.items.length > 0) {\n  inWhichCaseResults.items[0].insertText(\n    \" the actual implementation shall already be tested.\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 4) \"the cases I could elicit from the case description.\" ->\n//    \"the behavioral cases I could elicit from the case description,\n//    like \\u201can issue is allowed to transition to any issue state\\u201d.\"\n// ---------------------------------------------------------------------\nconst casesResults = body.search(\"cases I could\", { matchCase: true });\ncasesResults.load(\"text\");\nawait context.sync();\n\nif (casesResults.items.length > 0) {\n  casesResults.items[0].insertText(\"behavioral \", Word.InsertLocation.before);\n  await context.sync();\n}\n\nconst elicitResults = body.search(\"elicit from the case description.\", { matchCase: true });\nelicitResults.load(\"text\");\nawait context.sync();\n\nif (elicitResults.items.length > 0) {\n  elicitResults.items[0].insertText(\n    \"elicit from the case description, like \\u201can issue is allowed to transition to any issue state\\u201d.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) \"happy-test\" -> \"happy path test\"\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.Text = \"happy-test\"\n$find.Replacement.Text = \"happy path test\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# 2) Insert a new paragraph right after the \"I added one happy path\n#    test ... case description.\" paragraph, describing DB seeding.\n# ---------------------------------------------------------------------\n$para = $d.Paragraphs.Item(5)\n$para.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item(6)\n$newPara.Range.Text = \"To seed the \u201cdatabase\u201d I used the service layer calls, normally I would Mock or Fake this layer to seed data necessary for each test. I skipped that due to time restraint.\"\n\n# ---------------------------------------------------------------------\n# 3) \"...in which case\" -> \"...in which case the actual implementation\n#    shall already be tested.\"\n# ---------------------------------------------------------------------\n$caseRange = $d.Content\n$caseRange.Find.Execute(\"in which case\") | Out-Null\n$caseRange.Collapse(0)  # wdCollapseEnd\n$caseRange.InsertAfter(\" the actual implementation shall already be tested.\")\n\n# ---------------------------------------------------------------------\n# 4) \"the cases I could elicit from the case description.\" ->\n#    \"the behavioral cases I could elicit from the case description,\n#    like \u201can issue is allowed to transition to any issue state\u201d.\"\n# ---------------------------------------------------------------------\n$casesRange = $d.Content\n$casesRange.Find.Execute(\"cases I could\") | Out-Null\n$casesRange.Collapse(1)  # wdCollapseStart\n$casesRange.InsertBefore(\"behavioral \")\n\n$find2 = $d.Content.Find\n$find2.Text = \"elicit from the case description.\"\n$replacement = \"elicit from the case description, like \u201can issue is allowed to transition to any issue state\u201d.\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null\n"}
